$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2016528925619835
$ws.Range("C2").Value = 0.5586776859504132
$ws.Range("J2").Value = 0.006611570247933884
$ws.Range("P2").Value = 0.1570247933884298
$ws.Range("S2").Value = 0.07603305785123966
$ws.Range("B3").Value = 0.002816901408450704
$ws.Range("C3").Value = 0.02535211267605634
$ws.Range("J3").Value = 0.02253521126760564
$ws.Range("P3").Value = 0.7633802816901408
$ws.Range("S3").Value = 0.1859154929577465
$ws.Range("J4").Value = 0.05319148936170213
$ws.Range("P4").Value = 0.6595744680851063
$ws.Range("S4").Value = 0.2872340425531915
$ws.Range("O5").Value = 0.1666666666666667
$ws.Range("P5").Value = 0.8333333333333334
$ws.Range("B6").Value = 0.05806451612903226
$ws.Range("D6").Value = 0.01720430107526882
$ws.Range("F6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.232258064516129
$ws.Range("O6").Value = 0.01075268817204301
$ws.Range("Q6").Value = 0.1483870967741935
$ws.Range("R6").Value = 0.06666666666666667
$ws.Range("S6").Value = 0.4
$ws.Range("B7").Value = 0.1153846153846154
$ws.Range("D7").Value = 0.01538461538461539
$ws.Range("F7").Value = 0.05384615384615385
$ws.Range("J7").Value = 0.1179487179487179
$ws.Range("O7").Value = 0.01282051282051282
$ws.Range("Q7").Value = 0.1769230769230769
$ws.Range("R7").Value = 0.08974358974358974
$ws.Range("S7").Value = 0.417948717948718
$ws.Range("B8").Value = 0.09497816593886463
$ws.Range("D8").Value = 0.0240174672489083
$ws.Range("E8").Value = 0.003275109170305677
$ws.Range("F8").Value = 0.05895196506550218
$ws.Range("J8").Value = 0.1157205240174673
$ws.Range("O8").Value = 0.01310043668122271
$ws.Range("Q8").Value = 0.1681222707423581
$ws.Range("R8").Value = 0.09497816593886463
$ws.Range("S8").Value = 0.4268558951965066
$ws.Range("B9").Value = 0.1023017902813299
$ws.Range("D9").Value = 0.01790281329923274
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.1253196930946291
$ws.Range("O9").Value = 0.005115089514066497
$ws.Range("Q9").Value = 0.1611253196930946
$ws.Range("R9").Value = 0.09462915601023018
$ws.Range("S9").Value = 0.4347826086956522
$ws.Range("B10").Value = 0.1078582434514638
$ws.Range("D10").Value = 0.02157164869029276
$ws.Range("E10").Value = 0.001155624036979969
$ws.Range("F10").Value = 0.0724191063174114
$ws.Range("J10").Value = 0.1332819722650231
$ws.Range("O10").Value = 0.01232665639445301
$ws.Range("Q10").Value = 0.2068567026194145
$ws.Range("R10").Value = 0.09322033898305085
$ws.Range("S10").Value = 0.3513097072419106
$ws.Range("G11").Value = 0.1319073083778966
$ws.Range("J11").Value = 0.09090909090909091
$ws.Range("K11").Value = 0.1818181818181818
$ws.Range("L11").Value = 0.5828877005347594
$ws.Range("S11").Value = 0.0124777183600713
$ws.Range("F12").Value = 0.002932551319648094
$ws.Range("G12").Value = 0.7800586510263929
$ws.Range("J12").Value = 0.1524926686217009
$ws.Range("K12").Value = 0.008797653958944282
$ws.Range("L12").Value = 0.02639296187683285
$ws.Range("S12").Value = 0.02932551319648094
$ws.Range("F13").Value = 0.01162790697674419
$ws.Range("G13").Value = 0.7558139534883721
$ws.Range("J13").Value = 0.2093023255813954
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("F15").Value = 0.02262443438914027
$ws.Range("H15").Value = 0.2217194570135747
$ws.Range("I15").Value = 0.06108597285067873
$ws.Range("J15").Value = 0.3642533936651584
$ws.Range("K15").Value = 0.06787330316742081
$ws.Range("M15").Value = 0.006787330316742082
$ws.Range("O15").Value = 0.08144796380090498
$ws.Range("S15").Value = 0.1742081447963801
$ws.Range("F16").Value = 0.0194647201946472
$ws.Range("H16").Value = 0.2043795620437956
$ws.Range("I16").Value = 0.0681265206812652
$ws.Range("J16").Value = 0.4160583941605839
$ws.Range("K16").Value = 0.09245742092457421
$ws.Range("M16").Value = 0.0170316301703163
$ws.Range("N16").Value = 0.004866180048661801
$ws.Range("O16").Value = 0.0681265206812652
$ws.Range("S16").Value = 0.1094890510948905
$ws.Range("F17").Value = 0.02149321266968326
$ws.Range("H17").Value = 0.1798642533936652
$ws.Range("I17").Value = 0.09389140271493213
$ws.Range("J17").Value = 0.416289592760181
$ws.Range("K17").Value = 0.08936651583710407
$ws.Range("M17").Value = 0.02262443438914027
$ws.Range("N17").Value = 0.001131221719457014
$ws.Range("O17").Value = 0.07918552036199095
$ws.Range("S17").Value = 0.09615384615384616
$ws.Range("F18").Value = 0.0162037037037037
$ws.Range("H18").Value = 0.1759259259259259
$ws.Range("I18").Value = 0.1041666666666667
$ws.Range("J18").Value = 0.3912037037037037
$ws.Range("K18").Value = 0.1018518518518518
$ws.Range("M18").Value = 0.0162037037037037
$ws.Range("O18").Value = 0.0787037037037037
$ws.Range("S18").Value = 0.1157407407407407
$ws.Range("F19").Value = 0.01695599515542996
$ws.Range("H19").Value = 0.2030682276947921
$ws.Range("I19").Value = 0.0843762616067824
$ws.Range("J19").Value = 0.3859507468712152
$ws.Range("K19").Value = 0.1037545417844166
$ws.Range("M19").Value = 0.02220427937020589
$ws.Range("N19").Value = 0.0004037141703673799
$ws.Range("O19").Value = 0.06943883730318934
$ws.Range("S19").Value = 0.1138473960436011
